$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131117036
$ws.Range("B4").Value = 56748
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 205998
$ws.Range("F4").Value = "Nordfladdermus"
$ws.Range("G4").Value = "Eptesicus nilssonii"
$ws.Range("H4").Value = "(A.Keyserling & Blasius, 1839)"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "443"
$ws.Range("I4").Style = "Normal"
$ws.Range("N4").Value = "autobox med tidsexpansion"
$ws.Range("P4").Value = "Tryggaröd 9, Sk"
$ws.Range("Q4").Value = 437326
$ws.Range("R4").Value = 6227846
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Skåne"
$ws.Range("U4").Value = "Östra Göinge"
$ws.Range("V4").Value = "Skåne"
$ws.Range("W4").Value = "Gryt"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2025-07-08"
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2025-07-11"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AC4").Value = "Glänta, inslag av ädellöv."
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "David Alvunger"
$ws.Range("AX4").Value = "Johan Eklöf, Enviro Planning"

# Row 5
$ws.Range("A5").Value = 131116964
$ws.Range("B5").Value = 56762
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 100092
$ws.Range("F5").Value = "Större brunfladdermus"
$ws.Range("G5").Value = "Nyctalus noctula"
$ws.Range("H5").Value = "(Schreber, 1774)"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "42"
$ws.Range("I5").Style = "Normal"
$ws.Range("N5").Value = "autobox med tidsexpansion"
$ws.Range("P5").Value = "Tryggaröd 9, Sk"
$ws.Range("Q5").Value = 437326
$ws.Range("R5").Value = 6227846
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Skåne"
$ws.Range("U5").Value = "Östra Göinge"
$ws.Range("V5").Value = "Skåne"
$ws.Range("W5").Value = "Gryt"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2025-07-08"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2025-07-11"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AC5").Value = "Glänta, inslag av ädellöv."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "David Alvunger"
$ws.Range("AX5").Value = "Johan Eklöf, Enviro Planning"

# Row 6
$ws.Range("A6").Value = 131116934
$ws.Range("B6").Value = 56769
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 206002
$ws.Range("F6").Value = "Brunlångöra"
$ws.Range("G6").Value = "Plecotus auritus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("I6").Style = "Normal"
$ws.Range("N6").Value = "autobox med tidsexpansion"
$ws.Range("P6").Value = "Tryggaröd 9, Sk"
$ws.Range("Q6").Value = 437326
$ws.Range("R6").Value = 6227846
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Skåne"
$ws.Range("U6").Value = "Östra Göinge"
$ws.Range("V6").Value = "Skåne"
$ws.Range("W6").Value = "Gryt"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2025-07-08"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2025-07-11"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AC6").Value = "Glänta, inslag av ädellöv."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "David Alvunger"
$ws.Range("AX6").Value = "Johan Eklöf, Enviro Planning"

# Row 7
$ws.Range("A7").Value = 131117021
$ws.Range("B7").Value = 56746
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 100015
$ws.Range("F7").Value = "Barbastell"
$ws.Range("G7").Value = "Barbastella barbastellus"
$ws.Range("H7").Value = "(Schreber, 1774)"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("I7").Style = "Normal"
$ws.Range("N7").Value = "autobox med tidsexpansion"
$ws.Range("P7").Value = "Tryggaröd 9, Sk"
$ws.Range("Q7").Value = 437326
$ws.Range("R7").Value = 6227846
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Skåne"
$ws.Range("U7").Value = "Östra Göinge"
$ws.Range("V7").Value = "Skåne"
$ws.Range("W7").Value = "Gryt"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2025-07-08"
$ws.Range("Y7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2025-07-11"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AC7").Value = "Glänta, inslag av ädellöv."
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = "David Alvunger"
$ws.Range("AX7").Value = "Johan Eklöf, Enviro Planning"

